$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 34
$ws.Range("F3").Value = 65
$ws.Range("F4").Value = 0
$ws.Range("F5").Value = 8871
$ws.Range("F6").Value = 144
$ws.Range("F8").Value = 7100
$ws.Range("F10").Value = 5399
$ws.Range("F11").Value = 0
$ws.Range("F12").Value = 0
$ws.Range("F13").Value = 1094
$ws.Range("F14").Value = 396
$ws.Range("F15").Value = 402
$ws.Range("F16").Value = 31
$ws.Range("F17").Value = 560
$ws.Range("F19").Value = 270
$ws.Range("F20").Value = 145
$ws.Range("F21").Value = 203
$ws.Range("F24").Value = 9995
$ws.Range("F25").Value = 1896
$ws.Range("F27").Value = 48
$ws.Range("F29").Value = 2052
$ws.Range("F30").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("F33").Value = 1031
$ws.Range("F36").Value = 1363
$ws.Range("F40").Value = 0
$ws.Range("F42").Value = 0
$ws.Range("F44").Value = 1100
$ws.Range("F45").Value = 1073
$ws.Range("F46").Value = 0
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 0
$ws.Range("F6").Value = 3
$ws.Range("F7").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("F14").Value = 5
$ws.Range("F15").Value = 94
$ws.Range("F16").Value = 1
$ws.Range("F17").Value = 7
$ws.Range("F21").Value = 2
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 34
$ws.Range("F3").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("F6").Value = 8871
$ws.Range("F7").Value = 144
$ws.Range("F9").Value = 0
$ws.Range("F12").Value = 3
$ws.Range("F13").Value = 5399
$ws.Range("F15").Value = 6129
$ws.Range("F16").Value = 1094
$ws.Range("F18").Value = 402
$ws.Range("F23").Value = 145
$ws.Range("F24").Value = 0
$ws.Range("F27").Value = 0
$ws.Range("F28").Value = 9995
$ws.Range("F29").Value = 1896
$ws.Range("F30").Value = 1809
$ws.Range("F31").Value = 48
$ws.Range("F32").Value = 2052
$ws.Range("F37").Value = 2038
$ws.Range("F38").Value = 0
$ws.Range("F39").Value = 1363
$ws.Range("F41").Value = 1194
$ws.Range("F42").Value = 0
$ws.Range("F43").Value = 102
$ws.Range("F45").Value = 1100
$ws.Range("F46").Value = 1073
$ws.Range("F47").Value = 975
$ws.Range("F48").Value = 0
